# Lithuania A Lyga - league base update (18-04-2024 00:36)
#
# The source feed re-ordered a handful of match rows (same match data,
# different row position). Net effect on the worksheet:
#   - rows 26 <-> 27            : swapped
#   - rows 100,101,102,103,104  : cyclically rotated
#       (100<-102, 101<-103, 102<-104, 103<-100, 104<-101)
#   - rows 136 <-> 137          : swapped
# Column A (the running rank number) and columns C/D/E (Div / Div Original
# Name / Date, identical for every row in the block) stay put; only B:AC
# (match id .. PL_AhUnder) travel with the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC (A is left untouched on purpose)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues([int]$row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value()
    }
    return $vals
}

function Set-RowValues([int]$row, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# --- snapshot every row that will move, BEFORE writing anything back -------
$affectedRows = @(26, 27, 100, 101, 102, 103, 104, 136, 137)
$snapshot = @{}
foreach ($r in $affectedRows) {
    $snapshot[$r] = Get-RowValues $r
}

# --- target mapping: new row <- old row -------------------------------------
$mapping = @{
    26  = 27
    27  = 26
    100 = 102
    101 = 103
    102 = 104
    103 = 100
    104 = 101
    136 = 137
    137 = 136
}

foreach ($destRow in $affectedRows) {
    $srcRow = $mapping[$destRow]
    Set-RowValues $destRow $snapshot[$srcRow]
}

Write-Output "Reordered rows 26/27, 100-104, 136/137"
